$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The published syllabus sheet was re-synced from the source database and a
# number of rows (13-23) were re-shuffled: a row was inserted for
# "Programa resumido:" / "Semestral", the old "5840535 - Messias Borges
# Silva" entry moved up next to "Objetivos:" (and was reused again next to
# "Método:"), and the trailing blank-label row was dropped, shrinking the
# sheet from 24 to 23 rows. Apply the resulting layout directly, cell by
# cell, since the row-by-row content shuffle does not correspond to a
# simple insert/delete of a contiguous block.
# ---------------------------------------------------------------------------

$longAvaliados = "serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso.A média da disciplina será a média aritmética das duas provas."
$requisito = "LOB1012 -  Estatística  (Requisito fraco)`n"

# Row 10 - B/C text changes (long objectives text -> Messias Borges Silva)
$ws.Range("B10").Value = "5840535 - Messias Borges Silva"
$ws.Range("C10").Value = "5840535 - Messias Borges Silva"

# Row 13 becomes "Programa resumido:" / "Semestral" / "Semestral" with a new
# 60pt row height (previously it only held B13/C13 = Messias Borges Silva,
# row height unset).
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14 becomes "Short syllabus:" only (drop its old B/C long text), keep
# ht=60.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# Row 15 becomes "Programa:" with B/C = "01/01/1996", row height 120
# (previously "Programa:" alone with ht=60).
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/1996"
$ws.Range("C15").Value = "01/01/1996"
$ws.Rows.Item(15).RowHeight = 120

# Row 16 becomes "Syllabus:" only (drop its old B/C long text), ht stays 120.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Rows.Item(16).RowHeight = 120

# Row 17 becomes "Avaliação:" only, no explicit row height (was "Syllabus:"
# with ht=120).
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).RowHeight = $ws.Rows.Item(1).RowHeight()

# Row 18 becomes "Método:" with B/C = "5840535 - Messias Borges Silva" and a
# new ht=60 (previously "Avaliação:" alone, unset height).
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840535 - Messias Borges Silva"
$ws.Range("C18").Value = "5840535 - Messias Borges Silva"
$ws.Rows.Item(18).RowHeight = 60

# Row 19 becomes "Critério:" (was "Método:"), B/C stay "2 provas escritas".
$ws.Range("A19").Value = "Critério:"

# Row 20 becomes "Norma de recuperação:" (was "Critério:"), B/C stay the
# long "serão avaliados..." text.
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21 becomes "Bibliografia:" (was "Norma de recuperação:"); B/C stay
# "uma provas escrita..."; row height grows from 60 to 120.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# Row 22 becomes "Requisitos:" only (drop the old long bibliography text in
# B/C), and loses its explicit 120pt row height.
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).RowHeight = $ws.Rows.Item(1).RowHeight()

# Row 23 (was "Requisitos:" alone) becomes the old row 24's content: blank
# A, B/C = the requirement text, ht=30.
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = $requisito
$ws.Range("C23").Value = $requisito
$ws.Rows.Item(23).RowHeight = 30

# The old row 24 (duplicate of the requirement text) is no longer needed;
# remove it entirely so the sheet ends at row 23.
$ws.Rows.Item(24).Delete()
